# Correcciones lógicas y gráficas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Corrige el orden de los productos en el pedido de la fila 2
#    (antes "1-Torta,2-Queque," -> ahora "2-Queque,1-Torta,")
$ws.Range("A2").Value = "2-Queque,1-Torta,"

# 2) Agrega el pedido faltante de Alberto Hurtado (nueva fila 9)
$ws.Range("A9").Value = "2-Torta,"
$ws.Range("B9").Value = "18-10-2019"
$ws.Range("C9").Value = "20-10-2019"
$ws.Range("D9").Value = "Alberto Hurtado"
$ws.Range("E9").Value = "albertitohurtado@gmail.com"
$ws.Range("F9").Value = "9-48485930"
$ws.Range("G9").Value = "Finalizado"
$ws.Range("H9").Value = 10000.0
$ws.Range("I9").Value = 6000.0
$ws.Range("J9").Value = 0.0
$ws.Range("K9").Value = 7.0
